$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/22/2024  Through  4/28/2024"

# --- Data table updates (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = "'0"

# Row 15
$ws.Range("F15").Value = "'0"
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -25
$ws.Range("N15").Value = -10

# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 28.571428571428
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = -28.125
$ws.Range("I16").Value = 87
$ws.Range("J16").Value = 93
$ws.Range("K16").Value = -6.451612903225
$ws.Range("L16").Value = 10.126582278481
$ws.Range("M16").Value = 12.987012987013
$ws.Range("N16").Value = -66.666666666666

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -57.142857142857
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -13.793103448275
$ws.Range("I17").Value = 137
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 14.166666666666
$ws.Range("L17").Value = 8.730158730158
$ws.Range("M17").Value = 7.874015748031
$ws.Range("N17").Value = 5.384615384615

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 11.764705882352
$ws.Range("I18").Value = 59
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = -10.60606060606
$ws.Range("L18").Value = 20.408163265306
$ws.Range("M18").Value = -23.376623376623
$ws.Range("N18").Value = -78.148148148148

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -17.142857142857
$ws.Range("I19").Value = 128
$ws.Range("J19").Value = 134
$ws.Range("K19").Value = -4.477611940298
$ws.Range("L19").Value = -1.538461538461
$ws.Range("M19").Value = 58.024691358024
$ws.Range("N19").Value = 17.43119266055

# Row 20
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 80
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -17.857142857142
$ws.Range("I20").Value = 60
$ws.Range("J20").Value = 117
$ws.Range("K20").Value = -48.717948717948
$ws.Range("L20").Value = -24.050632911392
$ws.Range("M20").Value = 160.869565217391
$ws.Range("N20").Value = -51.219512195122

# Row 21
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -15.78947368421
$ws.Range("F21").Value = 120
$ws.Range("G21").Value = 145
$ws.Range("H21").Value = -17.241379310344
$ws.Range("I21").Value = 482
$ws.Range("J21").Value = 544
$ws.Range("K21").Value = -11.397058823529
$ws.Range("L21").Value = 1.26050420168
$ws.Range("M21").Value = 23.273657289002
$ws.Range("N21").Value = -47.207009857612

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 8
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 33.333333333333
$ws.Range("L22").Value = 14.285714285714
$ws.Range("M22").Value = 33.333333333333

# Row 23
$ws.Range("C23").Value = "'0"

# Row 24
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 61.111111111111
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 66
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 246
$ws.Range("J24").Value = 246
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = -8.550185873605
$ws.Range("M24").Value = 55.696202531645

# Row 25
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -22.222222222222
$ws.Range("I25").Value = 64
$ws.Range("J25").Value = 74
$ws.Range("K25").Value = -13.513513513513
$ws.Range("L25").Value = -48.387096774193

# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = -16.279069767441
$ws.Range("I26").Value = 197
$ws.Range("J26").Value = 151
$ws.Range("K26").Value = 30.46357615894
$ws.Range("L26").Value = 27.096774193548
$ws.Range("M26").Value = 14.53488372093

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -85.714285714285
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = 23.529411764705
$ws.Range("L27").Value = 0

# Row 28
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 39
$ws.Range("J28").Value = 34
$ws.Range("K28").Value = 14.705882352941
$ws.Range("L28").Value = -4.878048780487

# Row 29
$ws.Range("F29").Value = 2
$ws.Range("I29").Value = 8
$ws.Range("K29").Value = 33.333333333333
$ws.Range("L29").Value = 33.333333333333
$ws.Range("M29").Value = 60
$ws.Range("N29").Value = -66.666666666666

# Row 30
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 7
$ws.Range("K30").Value = 40
$ws.Range("L30").Value = 40
$ws.Range("M30").Value = 40
$ws.Range("N30").Value = -69.565217391304
